# Refactor: results are now saved to (and re-read from) a specified output
# folder. Re-running the sliding-window evaluation against the relocated
# results produces a new "IPC PO" (predicted) series; DELTA, DELTA^2 and the
# TOTAL / MSE summary rows are recomputed from it below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated column C ("IPC PO") values for data rows 2-51.
$newC = @(
    28.84267039079275,29.77707265575273,30.05592842394218,30.41447840443012,29.92317523131234,29.88122374981871,
    29.70272247807543,29.52518912846242,30.16943008269027,30.31824559916423,30.39605650836829,31.10876771424305,
    30.74681721347879,31.59678776628279,31.11618306159049,31.55355985879254,31.73061639384736,32.23293641710836,
    31.74654464723245,32.13029497748891,31.68627752760299,32.65755974196735,32.70016874744525,32.27986009465576,
    33.72895952859652,32.61288667413642,32.9278496063616,33.47847737626033,33.80295516145731,34.46276511415459,
    34.45397758051843,35.06007278628962,35.09758621209212,35.43159945849413,35.69388388264429,36.06136272225614,
    36.72716087155838,38.14990031659628,38.48782065156681,38.72278329549006,39.0000256948165,40.01617802748755,
    40.35817340030852,40.4896900035861,40.96970212329555,41.99817348785336,41.3252132420874,41.1035967236443,
    41.79583149888605,42.29066429366999
)

$sumDelta = 0.0
$sumDeltaSq = 0.0

for ($i = 0; $i -lt $newC.Count; $i++) {
    $row = $i + 2

    $ro    = $ws.Cells.Item($row, 2).Value2   # column B, "IPC RO" (unchanged)
    $po    = [double]$newC[$i]                # column C, "IPC PO" (new)
    $delta = $po - $ro                        # column D, "DELTA"
    $deltaSq = $delta * $delta                # column E, "DELTA^2"

    $ws.Cells.Item($row, 3).Value = $po
    $ws.Cells.Item($row, 4).Value = $delta
    $ws.Cells.Item($row, 5).Value = $deltaSq

    $sumDelta   += $delta
    $sumDeltaSq += $deltaSq
}

# Row 52 ("TOTAL"): sum of DELTA in column C, sum of DELTA^2 in column E.
$ws.Cells.Item(52, 3).Value = $sumDelta
$ws.Cells.Item(52, 5).Value = $sumDeltaSq

# Row 53 ("MSE"): mean of DELTA^2.
$ws.Cells.Item(53, 5).Value = $sumDeltaSq / $newC.Count
